$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new quarter's error row ("2020-04-01") needs to be inserted at row 3, which
# means every existing data row from 3 downward shifts down by one row. Shift
# manually (bottom-up) instead of using Rows.Insert() so no new/duplicate
# style is introduced for the blank row that Insert() would otherwise create.
for ($r = 22; $r -ge 3; $r--) {
    $dst = $r + 1
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($dst, $c).Value2 = $ws.Cells.Item($r, $c).Value2
    }
}

# Row 23 is brand new (the sheet used to end at row 22), so its A cell has no
# label formatting yet. Copy the bold/centered/bordered look used by every
# other row label (column A) down onto it before writing its text.
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)

# Populate the now-vacated row 3 with the new quarter's label and values.
$ws.Range("A3").Value = "2020-04-01 00:00:00_diff"
$ws.Range("B3").Value = 8.189548673647696
$ws.Range("C3").Value = -8.974959151229303
$ws.Range("D3").Value = -1.01725560823065
$ws.Range("E3").Value = 0.6812891314769711
$ws.Range("F3").Value = -2.226143026305237
$ws.Range("G3").Value = 0.06705817846833073
$ws.Range("H3").Value = -0.3695255935427733
